$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (RandomForestRegressor) - values updated
$ws.Range("B3").Value = 3425949512112.291
$ws.Range("C3").Value = 5394536991875.498
$ws.Range("D3").Value = 399874626461374.6

# Row 4 - label and values updated
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 3255548409134.465
$ws.Range("C4").Value = 3866524938045.253
$ws.Range("D4").Value = 209156961562312.7

# Row 5 - label and values updated
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 713974481241393.9
$ws.Range("C5").Value = 1265807397928304
$ws.Range("D5").Value = 3730740617340234
